# Update Courses pt. 2
# Adds course description text (column D) for several courses on Sheet1
# and updates the sheet's current view (scroll position / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 50 -> Creative Writing (Fall/Spring)
$ws.Range("D50").Value = "This course provides students an opportunity for additional writing instruction beyond the standard English program. Experimentation with many forms of writing is encouraged, with an emphasis on poetry, short stories, plays, and all forms of descriptive writing."

# Row 52 -> Journalism I
$ws.Range("D52").Value = "In Journalism 1 students learn the history and language of journalism; develop communication skills in writing, designing and editing for a variety of media; understand news and the process of publication; and become more critical readers and viewers of mass media."

# Row 53 -> Photojournalism I
$ws.Range("D53").Value = "In Photojournalism 1, students learn the principles of interviewing, copywriting, photography, layout, and design. Students will use publishing software to create pages for the school's yearbook. This course requires some after school time. This course requires an application and approval by instructor."

# Row 54 -> Film Study (F/S)
$ws.Range("D54").Value = "This course involves the study of classic and award-winning films. Students critique both the artistic and technical merits of the films. Students will produce short films each semester."

# Row 56 -> Debate (Fall)
$ws.Range("D56").Value = "Students will research the current topics and write cases both affirming and negating the resolutions. Students will also take part in regular in-class mock debates. This class can be used as preparation time for participation in the Woodson Debate team. Participation in at least one evening or Saturday debate event is required."

# Row 57 -> Forensics (Spring)
$ws.Range("D57").Value = "Students will be familiarized with some of the major forensic events, including Original Oratory, Dramatic Interpretation, Impromptu speaking, and Extemporaneous speaking. This class can be used as prep time for participation in the WTW Forensics (Speech) team. Participation in at least one evening or Saturday event is required."

# Match the D column style used by neighboring description cells (vertical-center alignment)
$ws.Range("D50").Style = $ws.Range("D49").Style
$ws.Range("D52:D54").Style = $ws.Range("D49").Style
$ws.Range("D56:D57").Style = $ws.Range("D49").Style

# Update the view: scroll so row 31 is the top-left visible row, and move the active selection to D57
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("D57").Select()
